$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Extend the green highlighted formatting from row 8 down to rows 9-10 ---
# Copy the (green, left-aligned, numFmt General) format from C8 onto A9:D10
$ws.Range("C8").Copy()
$ws.Range("A9:D10").PasteSpecial(-4122)

# Column B (the "KW" week-number column) keeps right alignment within the new
# green style, producing the dedicated green/right-aligned style used by B9:B10.
$ws.Range("B9:B10").HorizontalAlignment = -4152

# --- 2. Update the "Milestones" (column D) text content ---
$ws.Range("D7").Value = "Trace-cmd und KernelShark funktionieren"
$ws.Range("D8").Value = "Isolate CPUs of Host"
$ws.Range("D9").Value = "Decrease latency "
$ws.Range("D10").Value = "Preempt_RT Kernel Patch"

# --- 3. New milestone cells D11, D13, D15 need to be created with the
#        non-green (style used by rows 11-21) formatting. Use C11 (an
#        already-empty cell sharing that exact formatting) as the source. ---
$ws.Range("C11").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("D11").Value = "Zeitmessungen und Vergleiche zwischen verschiedenen Virtualisierungsmöglichkeiten unter Ubuntu sollen abgeschlossen sein"
$ws.Range("D12").Value = "QEMU/WSL unter Windows lauffähig"
$ws.Range("D13").Value = "Zeitmessungen und Vergleiche zwischen verschiedenen Virualisierungsmöglichkeiten unter Windows sollen abgeschlossen sein"
$ws.Range("D14").Value = "Zeitmessungen abgeschlossen Konklusio und Aufarbeitung/Vergleich PreemptRT gegen Xenomai"
$ws.Range("D15").Value = "Dedizierte Ressourcenzuteilung unter Windows und Messung Verhalten (Core-Sperrung,…)"

# --- 4. Restore the final selection state as left by the editing session ---
$ws.Range("E11").Select()
